$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns (history, electives, cs) before the existing
# "general_college_subjects.arts" column (old column R), shifting
# everything from R onward three columns to the right.
$ws.Range("R1:T1").EntireColumn.Insert()

# New header labels for the inserted columns
$ws.Range("R1").Value = "general_college_subjects.history"
$ws.Range("S1").Value = "general_college_subjects.electives"
$ws.Range("T1").Value = "general_college_subjects.cs"

# New data values for the inserted columns (row 2 - duke university)
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").Value = 0

# Normalize casing of existing categorical text values in row 2
$ws.Range("D2").Value = "considered"
$ws.Range("E2").Value = "considered"
$ws.Range("F2").Value = "considered"
$ws.Range("G2").Value = "very important"
$ws.Range("H2").Value = "very important"
$ws.Range("I2").Value = "considered"
$ws.Range("J2").Value = "considered"
